$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Codepen experiment task - time cost, time spent, notes
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = "3hr45min"
$ws.Range("E4").Value = "DONE Add <h2> with anchors; DONE Style menu; DONE JS autopopulate menu feature; Scroll between sections; NEW FORK: Tentative: inject aside with JS, inject anchors by scrubbing <h2> text nodes"

# Column E widened and row 3 height re-wrapped as a result of the longer note text
$ws.Columns.Item(5).ColumnWidth = 42.833333333333336
$ws.Rows.Item(3).RowHeight = 43.5
